# Updates cryptos list values/ranks to match the latest scrape (GitHub Actions run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.453.65'
$ws.Range("D3").Value = '2.659.83'
$ws.Range("E3").Value = '  +2.11%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = "'607.08"
$ws.Range("E5").Value = '  +2.01%  '
$ws.Range("D6").Value = "'157.34"
$ws.Range("E6").Value = '  +3.93%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("D8").Value = "'0.594"
$ws.Range("E8").Value = '  +1.23%  '
$ws.Range("E9").Value = '  +7.15%  '
$ws.Range("D10").Value = "'0.401"
$ws.Range("E10").Value = '  +4.17%  '
$ws.Range("D11").Value = "'5.89"
$ws.Range("E11").Value = '  +3.64%  '
$ws.Range("D13").Value = "'29.26"
$ws.Range("E13").Value = '  +5.78%  '
$ws.Range("D14").Value = '3.130.43'
$ws.Range("E14").Value = '  +1.77%  '
$ws.Range("D15").Value = "'0.0000181"
$ws.Range("E15").Value = '  +15.53%  '
$ws.Range("D16").Value = '65.286.79'
$ws.Range("E16").Value = '  +2.86%  '
$ws.Range("D17").Value = '2.815.10'
$ws.Range("E17").Value = '  +9.07%  '
$ws.Range("D18").Value = "'12.66"
$ws.Range("E18").Value = '  +0.99%  '
$ws.Range("D19").Value = "'4.85"
$ws.Range("E19").Value = '  +2.19%  '
$ws.Range("D20").Value = "'355.87"
$ws.Range("E20").Value = '  +2.47%  '
$ws.Range("E21").Value = '  +5.74%  '
$ws.Range("E22").Value = '  +0.21%  '
$ws.Range("E23").Value = '  +1.36%  '
$ws.Range("D24").Value = "'1.74"
$ws.Range("E24").Value = '  +2.23%  '
$ws.Range("D25").Value = "'9.58"
$ws.Range("E25").Value = '  +2.45%  '
$ws.Range("D26").Value = "'1.67"
$ws.Range("E26").Value = '  -1.68%  '
$ws.Range("D27").Value = "'8.34"
$ws.Range("E27").Value = '  +3.22%  '
$ws.Range("E28").Value = '  +2.28%  '
$ws.Range("B29").Value = 'Bittensor'
$ws.Range("C29").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D29").Value = "'548.69"
$ws.Range("E29").Value = '  -3.74%  '
$ws.Range("B30").Value = 'PEPE'
$ws.Range("C30").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D30").Value = '0.0₃0957'
$ws.Range("E30").Value = '  +12.01%  '
$ws.Range("D31").Value = "'0.998"
$ws.Range("E31").Value = '  -0.47%  '
$ws.Range("D32").Value = "'2.13"
$ws.Range("E32").Value = '  +3.43%  '
$ws.Range("E33").Value = '  +3.25%  '
$ws.Range("D34").Value = "'5.81"
$ws.Range("E34").Value = '  +10.64%  '
$ws.Range("D35").Value = "'6.48"
$ws.Range("E35").Value = '  +4.98%  '
$ws.Range("D36").Value = "'0.431"
$ws.Range("E36").Value = '  +3.90%  '
$ws.Range("E37").Value = '  +5.76%  '
$ws.Range("D38").Value = "'165.63"
$ws.Range("D39").Value = "'20.28"
$ws.Range("E39").Value = '  +3.01%  '
$ws.Range("E41").Value = '  -0.02%  '
$ws.Range("B42").Value = 'Aave'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D42").Value = "'168.58"
$ws.Range("E42").Value = '  +0.07%  '
$ws.Range("B43").Value = 'OKB'
$ws.Range("C43").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D43").Value = "'42.36"
$ws.Range("E43").Value = '  +6.83%  '
$ws.Range("D44").Value = "'4.14"
$ws.Range("E44").Value = '  +4.31%  '
$ws.Range("D45").Value = "'0.0614"
$ws.Range("E45").Value = '  +4.24%  '
$ws.Range("D46").Value = "'23.44"
$ws.Range("E46").Value = '  +5.23%  '
$ws.Range("E47").Value = '  +10.36%  '
$ws.Range("D48").Value = "'0.658"
$ws.Range("E48").Value = '  +4.07%  '
$ws.Range("D49").Value = "'0.0255"
$ws.Range("E49").Value = '  +0.52%  '
$ws.Range("E50").Value = '  +2.03%  '
$ws.Range("D51").Value = "'19.73"
$ws.Range("E51").Value = '  +2.68%  '
